# IEC_104_Conf.xlsx — "Repository for Connections and Main Configuration"
#
# Changes applied (per the target diff):
#   1. On "Main Configuration": reword the three header cells so they read
#      with spaces ("ExecutionTimeDefault" -> "Execution Time Default", etc.)
#   2. Move the active selection on "Connections" to E1 (was E17) without
#      leaving that sheet as the active tab.
#   3. Make "Main Configuration" the active sheet/tab (it was "Objects"
#      before) and leave its selection on C8 (was H11 on Main Configuration,
#      and "Objects" loses tabSelected / the workbook's activeTab pointer).

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main Configuration")
$wsConn = $wb.Worksheets.Item("Connections")

# 1. Reword the header row on "Main Configuration".
$wsMain.Range("A1").Value = "Execution Time Default"
$wsMain.Range("B1").Value = "Execution Time Short"
$wsMain.Range("C1").Value = "Execution Time Long"

# 2. Update the cursor position remembered on "Connections" (E17 -> E1)
#    while it stays in the background (not the active tab in the result).
$wsConn.Range("E1").Select()

# 3. Activate "Main Configuration" last so it becomes the workbook's active
#    tab, and park the selection on C8 as recorded in the saved view state.
$wsMain.Activate()
$wsMain.Range("C8").Select()
